# Insert a period after the question/sub-question numbering (e.g. "1.1" -> "1.1.")
# throughout the document, matching the target diff exactly. Each paragraph's
# numbering run is rewritten with its own Find/Replace so that unrelated runs in
# that same paragraph (a trailing math equation, a quoted 'tails' run, etc.) are
# left completely untouched.
#
# A couple of paragraphs (2.1 / 2.2 / 2.3 / 4.3) hold the numbering text in one
# run that is immediately followed by a plain-space run with identical (empty)
# formatting. Any write into that paragraph's Range re-coalesces same-formatted
# sibling runs, which would fold the following run into the one we just edited.
# To stop that, after rewriting the numbering run we nudge the very next
# character's Bold flag on/off (a value-preserving round trip) which is enough
# to make the engine keep that following run distinct again.

$d = $word.ActiveDocument

function Get-ParagraphByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Set-NumberingRun($prefix, $old, $new, [bool]$splitNext) {
    $p = Get-ParagraphByPrefix $prefix
    $start = $p.Range.Start
    $r = $d.Range($start, $start + $old.Length)
    $r.Text = $new
    if ($splitNext) {
        $afterStart = $start + $new.Length
        $afterRange = $d.Range($afterStart, $afterStart + 1)
        $afterRange.Bold = 1
        $afterRange.Bold = 0
    }
}

# Q1 sub-questions - each is the sole run in its paragraph.
Set-NumberingRun "1.1 " "1.1 Find the mean of this data set." "1.1. Find the mean of this data set." $false
Set-NumberingRun "1.2 " "1.2 Find the median of this data set." "1.2. Find the median of this data set." $false
Set-NumberingRun "1.3 " "1.3 Why would it be important to find the median of this data set, instead of only finding the mean?" "1.3. Why would it be important to find the median of this data set, instead of only finding the mean?" $false
Set-NumberingRun "1.4 " "1.4 Why would it not be important to find the mode of this data set?" "1.4. Why would it not be important to find the mode of this data set?" $false
Set-NumberingRun "1.5 " "1.5 Can you think of any potential reasons for the existence of the outlier?" "1.5. Can you think of any potential reasons for the existence of the outlier?" $false

# Q2 sub-questions - numbering run is followed by a plain space run then an oMath.
Set-NumberingRun "2.1" "2.1" "2.1." $true
Set-NumberingRun "2.2" "2.2" "2.2." $true
Set-NumberingRun "2.3" "2.3" "2.3." $true

# Q3 sub-questions - each is the sole run in its paragraph.
Set-NumberingRun "3.1 " "3.1 A business wants to visualize how the budget allocations breakdown of each department compares to the totality of the budget." "3.1. A business wants to visualize how the budget allocations breakdown of each department compares to the totality of the budget." $false
Set-NumberingRun "3.2 " "3.2 A florist wants to visualize the distribution of flower stem lengths." "3.2. A florist wants to visualize the distribution of flower stem lengths." $false
Set-NumberingRun "3.3 " ("3.3 A teacher wants to visualize their students" + [char]0x2019 + " exam scores against the number of hours they studied.") ("3.3. A teacher wants to visualize their students" + [char]0x2019 + " exam scores against the number of hours they studied.") $false
Set-NumberingRun "3.4 " "3.4 A bakery wants to visualize the trends in its chocolate bread sales revenue over time." "3.4. A bakery wants to visualize the trends in its chocolate bread sales revenue over time." $false

# Q4 sub-questions.
Set-NumberingRun "4.1 " "4.1 A researcher wants to test whether there is a relationship between cholestrol levels and heart disease risk." "4.1. A researcher wants to test whether there is a relationship between cholestrol levels and heart disease risk." $false
Set-NumberingRun "4.2 " "4.2 A researcher wants to model the relationship between cholestrol levels and heart disease risk in terms of a linear function." "4.2. A researcher wants to model the relationship between cholestrol levels and heart disease risk in terms of a linear function." $false
# 4.3's numbering run is followed by a plain space run, then a '\u2018tails\u2019' run.
Set-NumberingRun "4.3 " "4.3 You flip a coin three times. Each flip has a 0.5 chance of the coin landing on" "4.3. You flip a coin three times. Each flip has a 0.5 chance of the coin landing on" $true
Set-NumberingRun "4.4 " "4.4 A biologist collects a random sample of 100 birds and calculates their mean wingspan. They want to find a range of values in which the mean wingspan of all birds lies, with 90% confidence." "4.4. A biologist collects a random sample of 100 birds and calculates their mean wingspan. They want to find a range of values in which the mean wingspan of all birds lies, with 90% confidence." $false
